$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("U2").Value = "{'max_features': 1, 'window_size': 226, 'n_estimator': 15}"
$ws.Range("V2").Value = 42.76856716792099

# Row 3
$ws.Range("S3").Value = "[863, 866, 934, 954, 1068]"
$ws.Range("U3").Value = "{'max_features': 2, 'window_size': 251, 'n_estimator': 45}"
$ws.Range("V3").Value = 54.92935482901521

# Row 4
$ws.Range("S4").Value = "[832, 852, 862, 868, 873]"
$ws.Range("U4").Value = "{'max_features': 3, 'window_size': 202, 'n_estimator': 29}"
$ws.Range("V4").Value = 58.2047609569272

# Row 5
$ws.Range("U5").Value = "{'max_features': 6, 'window_size': 200, 'n_estimator': 38}"
$ws.Range("V5").Value = 35.67835730698425

# Row 6
$ws.Range("S6").Value = "[963, 991]"
$ws.Range("U6").Value = "{'max_features': 2, 'window_size': 254, 'n_estimator': 17}"
$ws.Range("V6").Value = 55.26876066299155

# Row 7
$ws.Range("U7").Value = "{'max_features': 4, 'window_size': 200, 'n_estimator': 24}"
$ws.Range("V7").Value = 33.94053089502268

# Row 8
$ws.Range("S8").Value = "[1, 248, 268, 827]"
$ws.Range("T8").Value = 0
$ws.Range("U8").Value = "{'max_features': 3, 'window_size': 240, 'n_estimator': 23}"
$ws.Range("V8").Value = 61.38626400299836

# Row 9
$ws.Range("S9").Value = "[813, 892, 984, 1026]"
$ws.Range("T9").Value = 0.25
$ws.Range("U9").Value = "{'max_features': 6, 'window_size': 213, 'n_estimator': 30}"
$ws.Range("V9").Value = 55.47462417802308
